$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 2.367275367578336
$ws.Range("D2").Value = 3.983303952131626
$ws.Range("F2").Value = 3.186643161705301
$ws.Range("H2").Value = 2.315550409313988
$ws.Range("J2").Value = 1.552801859005456
$ws.Range("L2").Value = 0.7980774557353087

# Row 4 updates
$ws.Range("B4").Value = 2.76229894385294
$ws.Range("D4").Value = 3.859250682189979
$ws.Range("F4").Value = 3.186478240571118
$ws.Range("H4").Value = 1.588926508372374
$ws.Range("J4").Value = 2.414540344278287
